$d = $word.ActiveDocument

# --- Edit 1: rewrite the intro paragraph -------------------------------
# Merge the "een <spellchecked>vacuumklok</spellchecked> met" runs back
# into a single run (this also drops the <w:proofErr> spell-check markers
# that bracketed "vacuumklok").
$d.Content.Find.Execute("een vacuumklok met", $true, $false, $false, $false,
                         $false, $true, 1, $false, "een vacuumklok met", 2) | Out-Null

# Drop the trailing period after "toegevoegd" - the sentence continues now.
$d.Content.Find.Execute("toegevoegd.", $true, $false, $false, $false,
                         $false, $true, 1, $false, "toegevoegd", 2) | Out-Null

# Append the new continuation as its own run: split the paragraph in two,
# fill in the new text, then rejoin the paragraphs by deleting the mark
# between them. This keeps the new text in a dedicated <w:r>, matching a
# second, separately-typed sentence rather than merging it into the first
# run.
$introPara = $d.Paragraphs(2)
$introPara.Range.InsertParagraphAfter()
$d.Paragraphs(3).Range.Text = " (zoemer en wit licht maken). Maar misschien zijn die niet eens nodig."
$joinRange = $d.Range($d.Paragraphs(2).Range.End - 1, $d.Paragraphs(2).Range.End)
$joinRange.Delete()

# --- Edit 2: add "Geluid is beweging" as a new bullet under Vortex kanon
$vortexPara = $d.Paragraphs(4)
$vortexPara.Range.InsertParagraphAfter()
$d.Paragraphs(5).Range.Text = "Geluid is beweging"
